$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    334 = @{ F = 204655; G = 3452 }
    335 = @{ F = 129667; G = 2959 }
    336 = @{ F = 103121; G = 3269 }
    337 = @{ F = 106092; G = 3007 }
    338 = @{ F = 221753; G = 3110 }
    339 = @{ F = 651776; G = 5601 }
    340 = @{ F = 383591; G = 3326 }
    341 = @{ F = 294437; G = 3626 }
    342 = @{ F = 177993; G = 3027 }
    343 = @{ F = 130119; G = 2902 }
    344 = @{ F = 134424; G = 2495 }
    345 = @{ F = 285775; G = 3275 }
    346 = @{ F = 660403; G = 4762 }
    347 = @{ F = 333371; G = 2819 }
    348 = @{ F = 230971; G = 3205 }
    349 = @{ F = 155991; G = 2696 }
    350 = @{ F = 125449; G = 2737 }
    351 = @{ F = 145270; G = 2740 }
    352 = @{ F = 298105; G = 3464 }
    353 = @{ F = 698552; G = 5123 }
    354 = @{ F = 296606 }
    355 = @{ F = 217291; G = 3347 }
    356 = @{ F = 156826; G = 2817 }
    357 = @{ F = 133309; G = 2943 }
    358 = @{ F = 155742; G = 2725 }
    359 = @{ F = 309430; G = 3285 }
    360 = @{ F = 677826; G = 4608 }
    361 = @{ F = 301686; G = 2419 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
